# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45171 (2023-09-02) to 45172 (2023-09-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSerial = 45171
$newSerial = 45172

$lastRow = 299
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value = $newSerial
    }
}
